$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (Title, Description) before the old "NeedItem" column (D),
# shifting NeedItem/NeedCount/RewardItem/RewardCount from D:G to F:I.
$ws.Columns("D:E").Insert()

# Header row
$ws.Cells.Item(1, 4).Value = "Title"
$ws.Cells.Item(1, 5).Value = "Description"

# Data row
$ws.Cells.Item(2, 4).Value = "길버트에게 오래된 동전을 가져다주자."
$ws.Cells.Item(2, 5).Value = "길버트 아저씨가 도와달라 한다. \r\n 얻어먹은 것도 있으니 일단은 노력해보자."

# Column widths to roughly match the authored best-fit widths.
$ws.Columns("D").ColumnWidth = 35.714285714285715
$ws.Columns("E").ColumnWidth = 68.71428571428571
$ws.Columns("F").ColumnWidth = 9

# Update selection to match the new active cell.
$ws.Range("E2").Select()
